# Apply updated Season Ladder standings (re-ranked by Total Points, refined
# club-name matching, and participation points recalculated with a floor
# function) to rows 7-22 of the ladder sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Club, Participation Points, Performance Points,
#             Total Points, Adjusted Total Points, ICL Eligible Number
$rows = @(
    @(7,  "Hills Triathlon Club",              30, 54, 84, 84, 178),
    @(8,  "Cronulla Triathlon Club",            15, 65, 80, 80, 207),
    @(9,  "Balance Triathlon Club",             30, 43, 73, 73, 127),
    @(10, "Coogee Triathlon Club",              30, 40, 70, 70, 119),
    @(11, "Pulse Performance",                  30, 26, 56, 56, 31),
    @(12, "Engadine Triathlon Club",            45, 0,  45, 45, 8),
    @(13, "Concord Triathlon Club",             30, 15, 45, 45, 53),
    @(14, "BRAT Triathlon Club",                30, 9,  39, 39, 161),
    @(15, "Manly Vipers Triathlon Club",        30, 9,  39, 39, 33),
    @(16, "Macarthur Triathlon Club",           15, 19, 34, 34, 69),
    @(17, "FilOz Triathlon Club",               30, 0,  30, 30, 10),
    @(18, "Sydney South West Triathlon Club",   30, 0,  30, 30, 2),
    @(20, "Northern Suburbs Triathlon Club",    15, 9,  24, 24, 67),
    @(21, "Australian Chinese Triathlon Club",  0,  0,  0,  0,  3),
    @(22, "Hunters Hills Triathlon Club",       0,  0,  0,  0,  6)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
